$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dereference the existing strings used in the data rows so the shared-string pool
# can be rebuilt from scratch in the desired order.
$ws.Range("A2:D5").Value = "___tmp___"

# Recreate the six distinct label strings in the exact order required by the target file:
# ECs, MuSCs, Bdnf, Ntrk2, FAPs, Inflammatory-Mac
$ws.Range("Z1").Value = "ECs"
$ws.Range("Z2").Value = "MuSCs"
$ws.Range("Z3").Value = "Bdnf"
$ws.Range("Z4").Value = "Ntrk2"
$ws.Range("Z5").Value = "FAPs"
$ws.Range("Z6").Value = "Inflammatory-Mac"
$ws.Range("Z1:Z6").Value = ""

# Fill in the data rows (rows 2-9) with the updated TPM-based values

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bdnf"
$ws.Range("C2").Value = "Ntrk2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.08140533333333333
$ws.Range("H2").Value = 0.244216
$ws.Range("I2").Value = 0.1131514935296598
$ws.Range("J2").Value = 0.1131514935296598
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6189250000000001
$ws.Range("N2").Value = 1.856775
$ws.Range("O2").Value = 0.09614699503454774
$ws.Range("P2").Value = 0.09614699503454775
$ws.Range("Q2").Value = 0.05038379593333334
$ws.Range("R2").Value = 0.4534541634
$ws.Range("S2").Value = 0.01087917608654786
$ws.Range("T2").Value = 0.01087917608654786

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bdnf"
$ws.Range("C3").Value = "Ntrk2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.08140533333333333
$ws.Range("H3").Value = 0.244216
$ws.Range("I3").Value = 0.1131514935296598
$ws.Range("J3").Value = 0.1131514935296598
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.072121333333333
$ws.Range("N3").Value = 15.216364
$ws.Range("O3").Value = 0.7879294335349575
$ws.Range("P3").Value = 0.7879294335349576
$ws.Range("Q3").Value = 0.412897727847111
$ws.Range("R3").Value = 3.716079550623999
$ws.Range("S3").Value = 0.08915539220045923
$ws.Range("T3").Value = 0.08915539220045923

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bdnf"
$ws.Range("C4").Value = "Ntrk2"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.08140533333333333
$ws.Range("H4").Value = 0.244216
$ws.Range("I4").Value = 0.1131514935296598
$ws.Range("J4").Value = 0.1131514935296598
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.01220666666666667
$ws.Range("N4").Value = 0.03662
$ws.Range("O4").Value = 0.001896246426284896
$ws.Range("P4").Value = 0.001896246426284896
$ws.Range("Q4").Value = 0.0009936877688888888
$ws.Range("R4").Value = 0.008943189919999999
$ws.Range("S4").Value = 0.0002145631152344159
$ws.Range("T4").Value = 0.0002145631152344159

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Bdnf"
$ws.Range("C5").Value = "Ntrk2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.08140533333333333
$ws.Range("H5").Value = 0.244216
$ws.Range("I5").Value = 0.1131514935296598
$ws.Range("J5").Value = 0.1131514935296598
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7340256666666667
$ws.Range("N5").Value = 2.202077
$ws.Range("O5").Value = 0.1140273250042099
$ws.Range("P5").Value = 0.1140273250042099
$ws.Range("Q5").Value = 0.05975360407022222
$ws.Range("R5").Value = 0.537782436632
$ws.Range("S5").Value = 0.01290236212741827
$ws.Range("T5").Value = 0.01290236212741827

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Bdnf"
$ws.Range("C6").Value = "Ntrk2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6380313333333333
$ws.Range("H6").Value = 1.914094
$ws.Range("I6").Value = 0.8868485064703402
$ws.Range("J6").Value = 0.8868485064703401
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6189250000000001
$ws.Range("N6").Value = 1.856775
$ws.Range("O6").Value = 0.09614699503454774
$ws.Range("P6").Value = 0.09614699503454775
$ws.Range("Q6").Value = 0.3948935429833333
$ws.Range("R6").Value = 3.55404188685
$ws.Range("S6").Value = 0.08526781894799988
$ws.Range("T6").Value = 0.08526781894799988

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Bdnf"
$ws.Range("C7").Value = "Ntrk2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6380313333333333
$ws.Range("H7").Value = 1.914094
$ws.Range("I7").Value = 0.8868485064703402
$ws.Range("J7").Value = 0.8868485064703401
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.072121333333333
$ws.Range("N7").Value = 15.216364
$ws.Range("O7").Value = 0.7879294335349575
$ws.Range("P7").Value = 0.7879294335349576
$ws.Range("Q7").Value = 3.23617233713511
$ws.Range("R7").Value = 29.125551034216
$ws.Range("S7").Value = 0.6987740413344983
$ws.Range("T7").Value = 0.6987740413344983

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Bdnf"
$ws.Range("C8").Value = "Ntrk2"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6380313333333333
$ws.Range("H8").Value = 1.914094
$ws.Range("I8").Value = 0.8868485064703402
$ws.Range("J8").Value = 0.8868485064703401
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.01220666666666667
$ws.Range("N8").Value = 0.03662
$ws.Range("O8").Value = 0.001896246426284896
$ws.Range("P8").Value = 0.001896246426284896
$ws.Range("Q8").Value = 0.007788235808888888
$ws.Range("R8").Value = 0.07009412228
$ws.Range("S8").Value = 0.00168168331105048
$ws.Range("T8").Value = 0.00168168331105048

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Bdnf"
$ws.Range("C9").Value = "Ntrk2"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6380313333333333
$ws.Range("H9").Value = 1.914094
$ws.Range("I9").Value = 0.8868485064703402
$ws.Range("J9").Value = 0.8868485064703401
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7340256666666667
$ws.Range("N9").Value = 2.202077
$ws.Range("O9").Value = 0.1140273250042099
$ws.Range("P9").Value = 0.1140273250042099
$ws.Range("Q9").Value = 0.4683313748042222
$ws.Range("R9").Value = 4.214982373238
$ws.Range("S9").Value = 0.1011249628767916
$ws.Range("T9").Value = 0.1011249628767916